$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 12.70776290365992
$ws.Range("C2").Value = 12.69131016872478
$ws.Range("D2").Value = 4.182622127043061
$ws.Range("F2").Value = 20.78720747333976
$ws.Range("G2").Value = 22.84119725964578
$ws.Range("H2").Value = 12.71106418598867
$ws.Range("I2").Value = 20.07234005538666
$ws.Range("L2").Value = 11.03312828331078
$ws.Range("M2").Value = 14.17802749424162
$ws.Range("N2").Value = 17.28809807930063
$ws.Range("O2").Value = 18.56714811764315
$ws.Range("B3").Value = 12.20832658026121
$ws.Range("C3").Value = 12.57743168253566
$ws.Range("D3").Value = 4.108591130166841
$ws.Range("F3").Value = 20.7681487588767
$ws.Range("G3").Value = 22.78730416526058
$ws.Range("H3").Value = 12.74648852071505
$ws.Range("I3").Value = 20.16415224079659
$ws.Range("L3").Value = 11.04888789219697
$ws.Range("M3").Value = 14.08556281132667
$ws.Range("N3").Value = 17.32739653074592
$ws.Range("O3").Value = 18.60702234619806
$ws.Range("B4").Value = 11.89193475148971
$ws.Range("C4").Value = 12.5067366346511
$ws.Range("D4").Value = 4.061835208402136
$ws.Range("F4").Value = 20.76257606239422
$ws.Range("G4").Value = 22.76285955512033
$ws.Range("H4").Value = 12.77043138438719
$ws.Range("I4").Value = 20.22461299346757
$ws.Range("L4").Value = 11.06027358841669
$ws.Range("M4").Value = 14.03041156531913
$ws.Range("N4").Value = 17.35320943252564
$ws.Range("O4").Value = 18.63601129818921
$ws.Range("B5").Value = 11.76073905452347
$ws.Range("C5").Value = 12.47775028203376
$ws.Range("D5").Value = 4.042467594534772
$ws.Range("F5").Value = 20.76184924666492
$ws.Range("G5").Value = 22.75507969245062
$ws.Range("H5").Value = 12.78073920496492
$ws.Range("I5").Value = 20.25027833518541
$ws.Range("L5").Value = 11.06534335862155
$ws.Range("M5").Value = 14.00836291690018
$ws.Range("N5").Value = 17.36415244622912
$ws.Range("O5").Value = 18.64895525354677
$ws.Range("B6").Value = 11.73882322128641
$ws.Range("C6").Value = 12.47292689141063
$ws.Range("D6").Value = 4.039233026864216
$ws.Range("F6").Value = 20.76182187425209
$ws.Range("G6").Value = 22.75391976741738
$ws.Range("H6").Value = 12.78248407403709
$ws.Range("I6").Value = 20.25460205702457
$ws.Range("L6").Value = 11.06621116613827
$ws.Range("M6").Value = 14.00472799322111
$ws.Range("N6").Value = 17.36599515800569
$ws.Range("O6").Value = 18.6511728056407
$ws.Range("B7").Value = 11.89017430528027
$ws.Range("C7").Value = 12.5063464110985
$ws.Range("D7").Value = 4.061575264774799
$ws.Range("F7").Value = 20.76256000561537
$ws.Range("G7").Value = 22.76274579276737
$ws.Range("H7").Value = 12.77056816911158
$ws.Range("I7").Value = 20.22495496760178
$ws.Range("L7").Value = 11.06034021980435
$ws.Range("M7").Value = 14.03011246170277
$ws.Range("N7").Value = 17.35335529592619
$ws.Range("O7").Value = 18.63618128998401
$ws.Range("B8").Value = 12.53768222413327
$ws.Range("C8").Value = 12.65221543866398
$ws.Range("D8").Value = 4.157373703895566
$ws.Range("F8").Value = 20.77936596575148
$ws.Range("G8").Value = 22.82082582156576
$ws.Range("H8").Value = 12.72282331894077
$ws.Range("I8").Value = 20.10314801636012
$ws.Range("L8").Value = 11.03820758744758
$ws.Range("M8").Value = 14.14581988828823
$ws.Range("N8").Value = 17.30129920621343
$ws.Range("O8").Value = 18.57996020247926
$ws.Range("B9").Value = 13.72312888668779
$ws.Range("C9").Value = 12.9313959827734
$ws.Range("D9").Value = 4.334381108428358
$ws.Range("F9").Value = 20.86077109379478
$ws.Range("G9").Value = 23.00285259179319
$ws.Range("H9").Value = 12.64660385267464
$ws.Range("I9").Value = 19.89675445547395
$ws.Range("L9").Value = 11.00835666687917
$ws.Range("M9").Value = 14.38478018275042
$ws.Range("N9").Value = 17.21254327631753
$ws.Range("O9").Value = 18.50555719838466
$ws.Range("B10").Value = 14.53457379173383
$ws.Range("C10").Value = 13.13135559138001
$ws.Range("D10").Value = 4.457142801667773
$ws.Range("F10").Value = 20.94977542106611
$ws.Range("G10").Value = 23.17728706770177
$ws.Range("H10").Value = 12.60123588367079
$ws.Range("I10").Value = 19.76496004842377
$ws.Range("L10").Value = 10.99466886271093
$ws.Range("M10").Value = 14.56657122816658
$ws.Range("N10").Value = 17.15541483117702
$ws.Range("O10").Value = 18.47285550969034
$ws.Range("B11").Value = 14.88943867919518
$ws.Range("C11").Value = 13.2210020943582
$ws.Range("D11").Value = 4.511279262628511
$ws.Range("F11").Value = 20.99650325521111
$ws.Range("G11").Value = 23.26523904279984
$ws.Range("H11").Value = 12.58290838826547
$ws.Range("I11").Value = 19.70932318769093
$ws.Range("L11").Value = 10.9902267666829
$ws.Range("M11").Value = 14.65037741689308
$ws.Range("N11").Value = 17.13117128143998
$ws.Range("O11").Value = 18.4627638620713
$ws.Range("B12").Value = 15.02167223785595
$ws.Range("C12").Value = 13.25474302439954
$ws.Range("D12").Value = 4.531523685942976
$ws.Range("F12").Value = 21.01508419394041
$ws.Range("G12").Value = 23.29975643853343
$ws.Range("H12").Value = 12.57630072712931
$ws.Range("I12").Value = 19.68887682617747
$ws.Range("L12").Value = 10.98880068698967
$ws.Range("M12").Value = 14.68225079300929
$ws.Range("N12").Value = 17.12224106836674
$ws.Range("O12").Value = 18.45963111409364
$ws.Range("B13").Value = 14.99329018841697
$ws.Range("C13").Value = 13.24748575153866
$ws.Range("D13").Value = 4.52717522303267
$ws.Range("F13").Value = 21.01104323698746
$ws.Range("G13").Value = 23.29226902519771
$ws.Range("H13").Value = 12.57770901000628
$ws.Range("I13").Value = 19.6932526261088
$ws.Range("L13").Value = 10.98909644107005
$ws.Range("M13").Value = 14.67538050848112
$ws.Range("N13").Value = 17.12415322644433
$ws.Range("O13").Value = 18.46027516805644
$ws.Range("B14").Value = 14.90036111905312
$ws.Range("C14").Value = 13.22378221074353
$ws.Range("D14").Value = 4.512949967862942
$ws.Range("F14").Value = 20.99801422941733
$ws.Range("G14").Value = 23.26805467366338
$ws.Range("H14").Value = 12.58235810537159
$ws.Range("I14").Value = 19.70762857754173
$ws.Range("L14").Value = 10.99010431428535
$ws.Range("M14").Value = 14.65299702148732
$ws.Range("N14").Value = 17.13043157434505
$ws.Range("O14").Value = 18.46249232283846
$ws.Range("B15").Value = 14.84315731737011
$ws.Range("C15").Value = 13.20923574861212
$ws.Range("D15").Value = 4.504202973505816
$ws.Range("F15").Value = 20.99014863274348
$ws.Range("G15").Value = 23.25337972381611
$ws.Range("H15").Value = 12.58524912829052
$ws.Range("I15").Value = 19.71651531968778
$ws.Range("L15").Value = 10.99075499203359
$ws.Range("M15").Value = 14.63930378995092
$ws.Range("N15").Value = 17.13430982289745
$ws.Range("O15").Value = 18.46394010346377
$ws.Range("B16").Value = 14.51108707116289
$ws.Range("C16").Value = 13.1254690368657
$ws.Range("D16").Value = 4.453569590823665
$ws.Range("F16").Value = 20.94684628247445
$ws.Range("G16").Value = 23.17171020314955
$ws.Range("H16").Value = 12.6024801507135
$ws.Range("I16").Value = 19.76868303086074
$ws.Range("L16").Value = 10.9949950270548
$ws.Range("M16").Value = 14.56111475906313
$ws.Range("N16").Value = 17.157034259194
$ws.Range("O16").Value = 18.47361137083675
$ws.Range("B17").Value = 14.30364723246644
$ws.Range("C17").Value = 13.07373157906143
$ws.Range("D17").Value = 4.422062605836294
$ws.Range("F17").Value = 20.92187228875226
$ws.Range("G17").Value = 23.12379462383211
$ws.Range("H17").Value = 12.61364285622319
$ws.Range("I17").Value = 19.80179286816051
$ws.Range("L17").Value = 10.99805283903564
$ws.Range("M17").Value = 14.5134169172086
$ws.Range("N17").Value = 17.17142136246748
$ws.Range("O17").Value = 18.48077045444373
$ws.Range("B18").Value = 14.18299533084551
$ws.Range("C18").Value = 13.04385067795543
$ws.Range("D18").Value = 4.403780382675728
$ws.Range("F18").Value = 20.90809593533315
$ws.Range("G18").Value = 23.09704606958049
$ws.Range("H18").Value = 12.62028082045236
$ws.Range("I18").Value = 19.82124303360698
$ws.Range("L18").Value = 10.99997962186756
$ws.Range("M18").Value = 14.48608802642946
$ws.Range("N18").Value = 17.17986067375099
$ws.Range("O18").Value = 18.48533849188179
$ws.Range("B19").Value = 14.14191797913313
$ws.Range("C19").Value = 13.03371293638706
$ws.Range("D19").Value = 4.397563136618069
$ws.Range("F19").Value = 20.90353280795248
$ws.Range("G19").Value = 23.08812951109152
$ws.Range("H19").Value = 12.62256566220378
$ws.Range("I19").Value = 19.82789826173835
$ws.Range("L19").Value = 11.00066086936246
$ws.Range("M19").Value = 14.47685373435124
$ws.Range("N19").Value = 17.18274630346622
$ws.Range("O19").Value = 18.48696246536672
$ws.Range("B20").Value = 14.32586876375347
$ws.Range("C20").Value = 13.07925196200382
$ws.Range("D20").Value = 4.4254332382847
$ws.Range("F20").Value = 20.92447003028871
$ws.Range("G20").Value = 23.12881153273427
$ws.Range("H20").Value = 12.61243205740461
$ws.Range("I20").Value = 19.7982262133181
$ws.Range("L20").Value = 10.99770994516962
$ws.Range("M20").Value = 14.51848365598054
$ws.Range("N20").Value = 17.1698728376644
$ws.Range("O20").Value = 18.47996174458903
$ws.Range("B21").Value = 14.92771557157753
$ws.Range("C21").Value = 13.23075024667185
$ws.Range("D21").Value = 4.517135292726413
$ws.Range("F21").Value = 21.00181721085599
$ws.Range("G21").Value = 23.27513433687901
$ws.Range("H21").Value = 12.58098352550835
$ws.Range("I21").Value = 19.70338911756792
$ws.Range("L21").Value = 10.98980133345883
$ws.Range("M21").Value = 14.65956802753269
$ws.Range("N21").Value = 17.12858068288496
$ws.Range("O21").Value = 18.46182239464338
$ws.Range("B22").Value = 15.30851351695131
$ws.Range("C22").Value = 13.32855329375342
$ws.Range("D22").Value = 4.575572150217563
$ws.Range("F22").Value = 21.05752695506977
$ws.Range("G22").Value = 23.37781519382322
$ws.Range("H22").Value = 12.5623686362246
$ws.Range("I22").Value = 19.64503436093638
$ws.Range("L22").Value = 10.98612470329616
$ws.Range("M22").Value = 14.75256786666709
$ws.Range("N22").Value = 17.10305252219388
$ws.Range("O22").Value = 18.45398184442754
$ws.Range("B23").Value = 15.10645013404494
$ws.Range("C23").Value = 13.27647018972868
$ws.Range("D23").Value = 4.544523362806483
$ws.Range("F23").Value = 21.02732559396424
$ws.Range("G23").Value = 23.32237627512529
$ws.Range("H23").Value = 12.57212628805301
$ws.Range("I23").Value = 19.67584709943421
$ws.Range("L23").Value = 10.98795067305931
$ws.Range("M23").Value = 14.70286671110003
$ws.Range("N23").Value = 17.1165440967457
$ws.Range("O23").Value = 18.45779902093464
$ws.Range("B24").Value = 14.31582673781545
$ws.Range("C24").Value = 13.07675662028608
$ws.Range("D24").Value = 4.423909899222409
$ws.Range("F24").Value = 20.92329377910987
$ws.Range("G24").Value = 23.12654089832796
$ws.Range("H24").Value = 12.61297877311464
$ws.Range("I24").Value = 19.79983740575198
$ws.Range("L24").Value = 10.99786444163549
$ws.Range("M24").Value = 14.51619269132034
$ws.Range("N24").Value = 17.17057240245661
$ws.Range("O24").Value = 18.48032595385224
$ws.Range("B25").Value = 13.41236759686366
$ws.Range("C25").Value = 12.85670614079146
$ws.Range("D25").Value = 4.287729364564778
$ws.Range("F25").Value = 20.83359145171703
$ws.Range("G25").Value = 22.94638993247635
$ws.Range("H25").Value = 12.66535758752067
$ws.Range("I25").Value = 19.94910839729458
$ws.Range("L25").Value = 11.01498288780364
$ws.Range("M25").Value = 14.31895848987504
$ws.Range("N25").Value = 17.23513200502069
$ws.Range("O25").Value = 18.52183492305214
